$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Hue step changes from 360/15 (15 buckets) to 360/16 (16 buckets) now that
#    an "if-else" node type has been added.
$ws.Range("D4").Formula = "=360/16"

# 2. The node-name list in column J is re-sorted (alphabetically, as driven by
#    the sortState) now that "if-else" is included. Re-assign every row's
#    label to match the new order.
$names = @("text","wildcard","characterset","anchor","whitespace","comment","or","concatenate","group","quantifier","reference","lookaround","if-else","flags","output")
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 4 + $i
    $ws.Range("J$row").Value = $names[$i]
}

# 3. The CSS variable-name formula pattern changes from
#    J4&" hsl(...)" to "--col-node-"&J4&": hsl(...)" (the "--col-node-" prefix
#    and colon move from the stored text in J into the formula itself).
$ws.Range("K4").Formula = '="--col-node-"&J4&": hsl("&ROUND(F4,0)&","&G4&", "&H4&");"'
$ws.Range("K5:K17").Formula = '="--col-node-"&J5&": hsl("&ROUND(F5,0)&","&G5&", "&H5&");"'

# 4. Rows 15-17 (lookaround / if-else / flags) now use a lightness of 70%
#    instead of 75%.
$ws.Range("H15").Value = "70%"
$ws.Range("H16").Value = "70%"
$ws.Range("H17").Value = "70%"

# 5. A new row 18 is populated for the 15th node ("output"), which used to be
#    the 9th entry and is now pushed to the end of the (growing) list.
$ws.Range("E18").Value = 15
$ws.Range("F18").Formula = "=(E18-1)*`$D`$4"
$ws.Range("G17").Copy($ws.Range("G18"))
$ws.Range("H14").Copy($ws.Range("H18"))
$ws.Range("J18").Value = "output"
$ws.Range("K18").Formula = '="--col-node-"&J18&": hsl("&ROUND(F18,0)&","&G18&", "&H18&");"'

# 6. Move the active selection to match where editing finished.
$ws.Range("L16").Select() | Out-Null
